# Deploy the implementation guide.
# Regenerates the "Metadata" sheet (refreshed Date / Contact, new Jurisdiction
# row) and renames the "Include from Ferlab.bio CodeS" sheet to "Include #0",
# matching a re-run of the FHIR IG publisher's xlsx exporter.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# --- Date: refresh the publication timestamp ---------------------------
$ws1.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# --- Contact: now resolves to a display string with url -----------------
$ws1.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# --- Insert a new "Jurisdiction" property row right after "Contact" ----
$ws1.Rows.Item(11).Insert()

# Match the formatting of the surrounding property rows (style-only copy).
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# --- Rename the code-system include sheet --------------------------------
$ws2 = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$ws2.Name = "Include #0"
